# Auto-generated Excel COM-interop script to apply the Balmung_Profits.xlsx cell updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 13833
$ws.Range("J3").Value = 13833
$ws.Range("L3").Value = 13833
$ws.Range("N3").Value = -14061
$ws.Range("H17").Value = 52772.95
$ws.Range("J17").Value = 55497.844
$ws.Range("L17").Value = 166493.532
$ws.Range("N17").Value = -166829.532
$ws.Range("H28").Value = 1038.6364
$ws.Range("I28").Value = 1157
$ws.Range("J28").Value = 506
$ws.Range("K28").Value = 1157
$ws.Range("L28").Value = 506
$ws.Range("M28").Value = -672
$ws.Range("N28").Value = -1476
$ws.Range("H74").Value = 12813.556
$ws.Range("I74").Value = 11626.286
$ws.Range("K74").Value = 11626.286
$ws.Range("M74").Value = -10690.286
$ws.Range("H77").Value = 12813.556
$ws.Range("I77").Value = 11626.286
$ws.Range("K77").Value = 58131.43
$ws.Range("M77").Value = -53451.43
$ws.Range("H102").Value = 13833
$ws.Range("J102").Value = 13833
$ws.Range("L102").Value = 13833
$ws.Range("N102").Value = -20323
$ws.Range("H137").Value = 9092618
$ws.Range("I137").Value = 1680.2
$ws.Range("K137").Value = 5040.6
$ws.Range("M137").Value = -2490.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 160373.61
$ws.Range("I32").Value = 167792.34
$ws.Range("J32").Value = 11998.667
$ws.Range("K32").Value = 167792.34
$ws.Range("L32").Value = 11998.667
$ws.Range("M32").Value = -167505.34
$ws.Range("N32").Value = -12572.667
$ws.Range("H44").Value = 99999
$ws.Range("J44").Value = 99999
$ws.Range("L44").Value = 99999
$ws.Range("N44").Value = -100975
$ws.Range("H45").Value = 47112
$ws.Range("I45").Value = 63997.938
$ws.Range("K45").Value = 63997.938
$ws.Range("M45").Value = -63620.938
$ws.Range("H61").Value = 1282317.6
$ws.Range("I61").Value = 3113.3823
$ws.Range("K61").Value = 3113.3823
$ws.Range("M61").Value = -2901.3823
$ws.Range("H74").Value = 399329.12
$ws.Range("I74").Value = 967.2143
$ws.Range("K74").Value = 967.2143
$ws.Range("M74").Value = -93.21429999999998
$ws.Range("H77").Value = 399329.12
$ws.Range("I77").Value = 967.2143
$ws.Range("K77").Value = 4836.0715
$ws.Range("M77").Value = -468.0715
$ws.Range("H136").Value = 1282317.6
$ws.Range("I136").Value = 3113.3823
$ws.Range("K136").Value = 9340.1469
$ws.Range("M136").Value = -6790.1469

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8604.684999999999
$ws.Range("I99").Value = 16999.25
$ws.Range("K99").Value = 16999.25
$ws.Range("M99").Value = -15501.25
$ws.Range("H107").Value = 9288.886
$ws.Range("I107").Value = 10202.333
$ws.Range("K107").Value = 10202.333
$ws.Range("M107").Value = -8282.333000000001
$ws.Range("H134").Value = 13434390
$ws.Range("I134").Value = 1420.6111
$ws.Range("J134").Value = 69232880
$ws.Range("K134").Value = 4261.8333
$ws.Range("L134").Value = 207698640
$ws.Range("M134").Value = -1726.8333
$ws.Range("N134").Value = -207703710

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1640.3334
$ws.Range("I16").Value = 1350.9048
$ws.Range("K16").Value = 1350.9048
$ws.Range("M16").Value = -1063.9048
$ws.Range("H22").Value = 606.2593000000001
$ws.Range("I22").Value = 433.23914
$ws.Range("J22").Value = 1601.125
$ws.Range("K22").Value = 433.23914
$ws.Range("L22").Value = 1601.125
$ws.Range("M22").Value = -83.23914000000002
$ws.Range("N22").Value = -2301.125
$ws.Range("H31").Value = 3651.2856
$ws.Range("I31").Value = 3512.8096
$ws.Range("J31").Value = 3755.1428
$ws.Range("K31").Value = 3512.8096
$ws.Range("L31").Value = 3755.1428
$ws.Range("M31").Value = -3217.8096
$ws.Range("N31").Value = -4345.1428
$ws.Range("H34").Value = 3651.2856
$ws.Range("I34").Value = 3512.8096
$ws.Range("J34").Value = 3755.1428
$ws.Range("K34").Value = 3512.8096
$ws.Range("L34").Value = 3755.1428
$ws.Range("M34").Value = -3310.8096
$ws.Range("N34").Value = -4159.1428
$ws.Range("H43").Value = 20612.25
$ws.Range("J43").Value = 20612.25
$ws.Range("L43").Value = 20612.25
$ws.Range("N43").Value = -20980.25
$ws.Range("H62").Value = 7202.615
$ws.Range("I62").Value = 5404.1113
$ws.Range("K62").Value = 5404.1113
$ws.Range("M62").Value = -4780.1113
$ws.Range("H65").Value = 7202.615
$ws.Range("I65").Value = 5404.1113
$ws.Range("K65").Value = 27020.5565
$ws.Range("M65").Value = -23900.5565
$ws.Range("H86").Value = 14445.632
$ws.Range("I86").Value = 9436.076999999999
$ws.Range("J86").Value = 25299.666
$ws.Range("K86").Value = 9436.076999999999
$ws.Range("L86").Value = 25299.666
$ws.Range("M86").Value = -8313.076999999999
$ws.Range("N86").Value = -27545.666
$ws.Range("H89").Value = 14445.632
$ws.Range("I89").Value = 9436.076999999999
$ws.Range("J89").Value = 25299.666
$ws.Range("K89").Value = 47180.38499999999
$ws.Range("L89").Value = 126498.33
$ws.Range("M89").Value = -41564.38499999999
$ws.Range("N89").Value = -137730.33
$ws.Range("H101").Value = 20612.25
$ws.Range("J101").Value = 20612.25
$ws.Range("L101").Value = 20612.25
$ws.Range("N101").Value = -27102.25
$ws.Range("H105").Value = 2645.818
$ws.Range("I105").Value = 2233.7778
$ws.Range("K105").Value = 2233.7778
$ws.Range("M105").Value = -486.7777999999998
$ws.Range("H113").Value = 1640.3334
$ws.Range("I113").Value = 1350.9048
$ws.Range("K113").Value = 1350.9048
$ws.Range("M113").Value = 819.0952
$ws.Range("H122").Value = 5928.1816
$ws.Range("I122").Value = 5776.25
$ws.Range("K122").Value = 17328.75
$ws.Range("M122").Value = -14878.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1812843.1
$ws.Range("I4").Value = 2097804.2
$ws.Range("J4").Value = 143785
$ws.Range("K4").Value = 6293412.600000001
$ws.Range("L4").Value = 431355
$ws.Range("M4").Value = -6293300.600000001
$ws.Range("N4").Value = -431579
$ws.Range("H129").Value = 2965.7144
$ws.Range("I129").Value = 583.75
$ws.Range("J129").Value = 6141.6665
$ws.Range("K129").Value = 1751.25
$ws.Range("L129").Value = 18424.9995
$ws.Range("M129").Value = 3248.75
$ws.Range("N129").Value = -28424.9995
$ws.Range("H131").Value = 4787392.5
$ws.Range("J131").Value = 2844.8462
$ws.Range("L131").Value = 8534.5386
$ws.Range("N131").Value = -18614.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H80").Value = 16824570
$ws.Range("I80").Value = 136170.38
$ws.Range("J80").Value = 83578170
$ws.Range("K80").Value = 136170.38
$ws.Range("L80").Value = 83578170
$ws.Range("M80").Value = -135172.38
$ws.Range("N80").Value = -83580166
$ws.Range("H83").Value = 16824570
$ws.Range("I83").Value = 136170.38
$ws.Range("J83").Value = 83578170
$ws.Range("K83").Value = 680851.9
$ws.Range("L83").Value = 417890850
$ws.Range("M83").Value = -675859.9
$ws.Range("N83").Value = -417900834
$ws.Range("H102").Value = 45456080
$ws.Range("I102").Value = 50001132
$ws.Range("K102").Value = 50001132
$ws.Range("M102").Value = -49999510
$ws.Range("H123").Value = 57942.57
$ws.Range("J123").Value = 57942.57
$ws.Range("L123").Value = 57942.57
$ws.Range("N123").Value = -62842.57
$ws.Range("H126").Value = 3442.2
$ws.Range("I126").Value = 3302.75
$ws.Range("K126").Value = 9908.25
$ws.Range("M126").Value = -7438.25
$ws.Range("H132").Value = 646633.2
$ws.Range("I132").Value = 6171.6665
$ws.Range("K132").Value = 18514.9995
$ws.Range("M132").Value = -15984.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 14611
$ws.Range("J2").Value = 14611
$ws.Range("L2").Value = 14611
$ws.Range("N2").Value = -14835
$ws.Range("H16").Value = 5252
$ws.Range("I16").Value = 5252
$ws.Range("K16").Value = 5252
$ws.Range("M16").Value = -5082
$ws.Range("H23").Value = 721085
$ws.Range("I23").Value = 721085
$ws.Range("K23").Value = 721085
$ws.Range("M23").Value = -720855
$ws.Range("H46").Value = 13749.8
$ws.Range("I46").Value = 15437.25
$ws.Range("K46").Value = 15437.25
$ws.Range("M46").Value = -15249.25
$ws.Range("H82").Value = 3180.4
$ws.Range("I82").Value = 3001
$ws.Range("J82").Value = 3300
$ws.Range("K82").Value = 3001
$ws.Range("L82").Value = 3300
$ws.Range("M82").Value = -2640
$ws.Range("N82").Value = -4022
$ws.Range("H85").Value = 3180.4
$ws.Range("I85").Value = 3001
$ws.Range("J85").Value = 3300
$ws.Range("K85").Value = 3001
$ws.Range("L85").Value = 3300
$ws.Range("M85").Value = -1753
$ws.Range("N85").Value = -5796
$ws.Range("H100").Value = 4529
$ws.Range("I100").Value = 3121.4285
$ws.Range("J100").Value = 6499.6
$ws.Range("K100").Value = 3121.4285
$ws.Range("L100").Value = 6499.6
$ws.Range("M100").Value = -2580.4285
$ws.Range("N100").Value = -7581.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 82449.5
$ws.Range("J98").Value = 82449.5
$ws.Range("L98").Value = 82449.5
$ws.Range("N98").Value = -88439.5
$ws.Range("H123").Value = 100469.664
$ws.Range("J123").Value = 100469.664
$ws.Range("L123").Value = 100469.664
$ws.Range("N123").Value = -110269.664
$ws.Range("H136").Value = 22632.914
$ws.Range("I136").Value = 26033.1
$ws.Range("K136").Value = 78099.29999999999
$ws.Range("M136").Value = -75549.29999999999
